# criando comandos para criar a chave ssh
# Adds the new "git pull" row plus the ssh-key creation instructions
# at the bottom of the worksheet (rows 49-55), matching the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49: git pull / traz a ultima versão no servidor remoto
$ws.Cells.Item(49, 1).Value = "git pull"
$ws.Cells.Item(49, 2).Value = "traz a ultima versão no servidor remoto"

# Row 50: ssh-keygen
$ws.Cells.Item(50, 1).Value = "ssh-keygen"

# Row 51: cd ~/.ssh/
$ws.Cells.Item(51, 1).Value = "cd ~/.ssh/"

# Row 52: start .
$ws.Cells.Item(52, 1).Value = "start ."

# Row 53: muda o nome dos arquivos
$ws.Cells.Item(53, 1).Value = "muda o nome dos arquivos"

# Row 54: eval $(ssh-agent)
$ws.Cells.Item(54, 1).Value = 'eval $(ssh-agent)'

# Row 55: ssh-add ~/.ssh/nome_do_arquivo (não tem extenção)
$ws.Cells.Item(55, 1).Value = "ssh-add ~/.ssh/nome_do_arquivo (não tem extenção)"

# Widen column A slightly to fit the new, longer entries.
$ws.Columns.Item(1).ColumnWidth = 47.45

# Move the on-screen selection down to where the new rows were typed,
# mirroring the author's cursor position after the edit.
$ws.Cells.Item(52, 2).Select()
